$d = $word.ActiveDocument

# 1) Remove the _GoBack bookmark from its current location (top of the
#    document, inside the very first paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Delete the whole paragraph that reads "Leer nombre completo, ..."
#    (it is being removed outright in the revision) and 3) re-create the
#    _GoBack bookmark at the very start of the paragraph that follows it
#    ("Imprimir nombre completo, ...").
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -like "Leer nombre completo, matrícula, carrera, escuela de procedencia y descripción general.*") {
        $p.Range.Delete()
        break
    }
}

$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -like "Imprimir nombre completo, matrícula, carrera, escuela de procedencia y descripción general.*") {
        $startPoint = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $startPoint)
        break
    }
}
